# Update cryptos list with the latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.507.23'
$ws.Range("E2").Value = '  -2.37%  '
$ws.Range("D3").Value = '2.331.31'
$ws.Range("E3").Value = '  -3.42%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.635'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.612'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.13%  '
$ws.Range("E11").Value = '  -3.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.89%  '
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.980'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -9.06%  '
$ws.Range("D16").Value = '2.681.77'
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '2.303.70'
$ws.Range("E17").Value = '  -4.53%  '
$ws.Range("D18").Value = '42.464.71'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.15%  '
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0885'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  -1.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.33%  '
$ws.Range("E37").Value = '  -8.57%  '
$ws.Range("E38").Value = '  -5.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.52%  '
$ws.Range("E40").Value = '  -11.64%  '
$ws.Range("E41").Value = '  -10.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.230'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '114.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.99%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.58%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.98%  '
$ws.Range("E48").Value = '  -4.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0993'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.74%  '
